$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.481.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.930.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.64'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000224'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.20'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.415.37'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.448.33'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.931.67'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '433.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.39'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.72'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.78%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.90'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.73'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.109'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0875'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.19%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.61'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.98'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.99'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.49'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.96'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.34%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0344'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.694.50'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '133.13'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '363.50'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.53'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.85%  '
